$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "2025-10-25"
$ws.Range("B71").Value = 53.81999969482422
$ws.Range("C71").Value = 403.2999877929688
$ws.Range("D71").Value = 326.6000061035156
